$d = $word.ActiveDocument

# The document currently ends with an empty "ListParagraph" item (the
# last paragraph in the body, just before the sectPr). We want to add a
# new list item "Add the Biggest method." right before that trailing
# empty paragraph, inheriting its ListParagraph style / numbering.

$d.Paragraphs.Last.Range.InsertParagraphBefore()

# Re-fetch paragraphs by (live) index after the mutation above - cached
# paragraph/range object references taken before the insertion can be
# stale, so look the new paragraph up fresh by position. The newly
# inserted paragraph is now the second-to-last paragraph in the body.
$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count - 1)
$newPara.Range.Text = "Add the Biggest method."
